# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the "_old" / "_new" column-header suffixes to the concrete
# format-version names ("_FV2304" / "_FV2310"), turns the data range
# A1:U68 into a real Excel Table ("Table1") with autofilter + banded
# rows, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1) -----------------------------
# Columns A:J carried the "_old" suffix -> "_FV2304";
# column K ("diff") is unchanged;
# columns L:U carried the "_new" suffix -> "_FV2310".
$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $fv2304Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2304Headers[$i]
}

for ($i = 0; $i -lt $fv2310Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2310Headers[$i]
}

# --- 2. Freeze the header row ----------------------------------------
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into a real Excel table -------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U68"), 0, 1)
$lo.Name = "Table1"
